$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing "9.5" text value (C2) to "9.55" -------------------
# The cell must stay a text value (it is read as the shared string "9.55",
# same as the other dimension cells in this row), not get converted into a
# number, and it must keep its original (default) style. Temporarily format
# as Text so the numeric-looking string isn't auto-coerced to a Double, then
# restore the original "Normal" style so no new cell format lingers behind.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "9.55"
$ws.Range("C2").Style = "Normal"

# --- Fill in the new dimension data for the L_SPT, M_SPT and S_SPT rows ----
$ws.Range("B3").Value = 2.4
$ws.Range("C3").Value = 4.8
$ws.Range("D3").Value = 0.4
$ws.Range("E3").Value = 1.2
$ws.Range("G3").Value = "rounded edges"

$ws.Range("B4").Value = 1.2
$ws.Range("C4").Value = 2.9
$ws.Range("D4").Value = 0.2
$ws.Range("E4").Value = 0.6

$ws.Range("B5").Value = 0.6
$ws.Range("C5").Value = 1.2
$ws.Range("D5").Value = 0.1
$ws.Range("E5").Value = 0.3

# --- Match the saved selection/view state -----------------------------------
$ws.Range("E7").Select()
